# Applies "Fixed initial nbr of workers in each firm" changes.
$wb = $excel.ActiveWorkbook

$wsInit = $wb.Worksheets.Item("Initialization_Parameters")
$wsMain = $wb.Worksheets.Item("Main_Loop_Parameters")

# --- Main_Loop_Parameters sheet updates ---
$wsMain.Range("B4").Value = 1500
$wsMain.Range("B5").Value = 100
$wsMain.Range("B6").Value = 30
$wsMain.Range("B7").Formula = "=B5+B6"

# --- Initialization_Parameters sheet updates ---
$wsInit.Range("B29").Value = 0
$wsInit.Range("B30").Value = 0
$wsInit.Range("B31").Value = 0
$wsInit.Range("B34").Value = 1
$wsInit.Range("B46").Formula = "=INT((Main_Loop_Parameters!B4-Main_Loop_Parameters!B7)/Main_Loop_Parameters!B7 )"
$wsInit.Range("B47").Value = 0

# --- Selection / view state ---
$wsInit.Range("B34").Select()
$wsMain.Range("B8").Select()

# Make Main_Loop_Parameters the active (selected) sheet/tab.
$wsMain.Activate()
